$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 2017-2022 data (columns C:N, rows 2-12) two columns to
# the right (E:P) to make room for the new 2016 columns (C:D). Columns are
# moved right-to-left so the source of each move is read before it gets
# overwritten by an earlier (more-left) move's destination.
for ($col = 14; $col -ge 3; $col--) {
    $destCol = $col + 2
    $src = $ws.Range($ws.Cells.Item(2, $col), $ws.Cells.Item(12, $col))
    $src.Cut()
    $dest = $ws.Cells.Item(2, $destCol)
    $dest.Select()
    $ws.Paste()
}

# ---- Header row 2 (inflow/otflow labels) for new 2016 columns ----
$ws.Range("C2").Value = "inflow"
$ws.Range("D2").Value = "otflow"

# ---- Year row 3 ----
$ws.Range("C3").Value = 2016
$ws.Range("D3").Value = 2016

# ---- Data rows 4-12 for year 2016 ----
$data2016 = @{
    4  = @(78934, 61398)
    5  = @(7219, 5449)
    6  = @(50606, 44840)
    7  = @(3638, 2601)
    8  = @(2299, 1331)
    9  = @(6639, 3523)
    10 = @(2342, 1097)
    11 = @(3964, 1662)
    12 = @(2227, 895)
}

foreach ($row in $data2016.Keys) {
    $vals = $data2016[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}

# Re-assert formatting across the whole table (rows 2-3 centered+bold,
# data rows 4-12 centered) since Cut/Paste into previously-empty cells
# (columns O:P) does not always retain direct formatting.
$ws.Range("C2:P3").HorizontalAlignment = -4108
$ws.Range("C2:P3").Font.Bold = $true
$ws.Range("C4:P12").HorizontalAlignment = -4108
$ws.Range("C4:P12").Font.Bold = $false

# Update view: selection + zoom
$ws.Range("E18").Select()
$excel.ActiveWindow.Zoom = 100
